# TE-Backend.xlsx: Add activity api & update docu
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix existing style (cellXfs index 11, currently only used by D16) ---
# Target: applyFont=true, wrapText=true, horizontal=left, vertical=top (unchanged)
$ws.Range("D16").WrapText = $true

# --- 2. Populate the new cell values, in the same left-to-right / top-to-bottom
#        order the workbook's shared-string table records them in ---
$ws.Range("A15").Value = "ACTIVITY"
$ws.Range("B15").Value = "List of Activity"
$ws.Range("C15").Value = "/api/activity/viewset/list/"
$ws.Range("D15").Value = "{`n    ""activity_name"": ""Activity #1"",`n    ""chapter"": 3`n}"

$ws.Range("C16").Value = "/api/activity/viewset/list/<int:pk>/"
$ws.Range("E16").Value = "get specific activity"

$ws.Range("B17").Value = "Set prof activity date"
$ws.Range("C17").Value = "/api/activity/viewset/prof_activity/"
$ws.Range("D17").Value = "{`n    ""activity"": 1,`n    ""section"": 5,`n    ""start"": ""2020-12-16"",`n    ""end"": ""2020-12-20"",`n    ""remarks"": false`n}"

$ws.Range("B18").Value = "List of Activity (Prof)"
$ws.Range("C18").Value = "/api/activity/viewset/prof_activity/"
$ws.Range("D18").Value = "{`n        ""id"": 1,`n        ""activity"": 1,`n        ""activity_name"": ""Activity #1"",`n        ""section"": 5,`n        ""section_code"": ""CEIT--03-1001E"",`n        ""start"": ""2020-11-16"",`n        ""end"": ""2020-11-20"",`n        ""remarks"": false`n    }"

$ws.Range("C19").Value = "/api/activity/viewset/prof_activity/<int:pk>/"
$ws.Range("E19").Value = "get specific activity (prof)"

# --- 3. Copy cell formatting from existing rows that already carry the
#        styles the new rows need, so the style table gets reused/deduped
#        rather than growing new near-duplicate entries ---
$ws.Range("B2").Copy()
$ws.Range("B15").PasteSpecial(-4122)

$ws.Range("D16").Copy()
$ws.Range("D15").PasteSpecial(-4122)

$ws.Range("B2").Copy()
$ws.Range("B17").PasteSpecial(-4122)

$ws.Range("C6").Copy()
$ws.Range("C17").PasteSpecial(-4122)

$ws.Range("D3").Copy()
$ws.Range("D17").PasteSpecial(-4122)

$ws.Range("B2").Copy()
$ws.Range("B18").PasteSpecial(-4122)

$ws.Range("B2").Copy()
$ws.Range("C18").PasteSpecial(-4122)

$ws.Range("D3").Copy()
$ws.Range("D18").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- 4. Row heights for the new rows ---
$ws.Rows.Item(15).RowHeight = 40.25
$ws.Rows.Item(16).RowHeight = 21.6
$ws.Rows.Item(17).RowHeight = 79.85
$ws.Rows.Item(18).RowHeight = 113.4
$ws.Rows.Item(19).RowHeight = 12.8

# --- 5. Merge the activity JSON cell across its two rows ---
$ws.Range("D15:D16").Merge()

# --- 6. Widen column D slightly to fit the new content ---
$ws.Columns.Item(4).ColumnWidth = 32.1

# --- 7. Park the view/selection near the newly added rows ---
$ws.Range("C14").Select()
